$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.040.34'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '1.864.36'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.63'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5133'
$ws.Range('E7').Value = '  +2.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3896'
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08336'
$ws.Range('E9').Value = '  -2.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.109'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.162'
$ws.Range('E11').Value = '  -2.35%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.862.81'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.43'
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.263'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.006'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001097'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.68'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06629'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.62'
$ws.Range('E19').Value = '  -3.05%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.987'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '28.082.60'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.02'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.248'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('B25').Value = 'LEO'
$ws.Range('C25').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.358'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').Value = '2.075.38'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.466'
$ws.Range('E27').Value = '  -5.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.25'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.48'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.48'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1060'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.029'
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.813'
$ws.Range('E33').Value = '  +2.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.601'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.432'
$ws.Range('E35').Value = '  -2.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02424'
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06515'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2175'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6474'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.986'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.213'
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.30'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6062'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.07'
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.288'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.675'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.994'
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.217'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.50'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.95'
$ws.Range('E51').Value = '  -3.69%  '
